$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Scroll the view so column B is the left-most visible column
$ws.Application.ActiveWindow.ScrollColumn = 2

# Widen column N slightly, and add a new column O
$ws.Columns.Item(14).ColumnWidth = 29.57421875
$ws.Columns.Item(15).ColumnWidth = 17.8515625

# New header for column O
$ws.Range("O1").Value = "must thermalize"

# New "yes"/"no" values for column O on the relevant rows
$ws.Range("O3").Value = "yes"
$ws.Range("O4").Value = "yes"
$ws.Range("O5").Value = "yes"
$ws.Range("O6").Value = "yes"
$ws.Range("O9").Value = "yes"
$ws.Range("O10").Value = "yes"
$ws.Range("O11").Value = "yes"
$ws.Range("O12").Value = "yes"
$ws.Range("O15").Value = "no"
$ws.Range("O18").Value = "yes"

# I10 changes from 0 to 1
$ws.Range("I10").Value = 1
